$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert 3 new price records for "Repollo" (Femacal de La Calera)
# right before the existing row 573, pushing the rest of the table down by 3 rows.
$ws.Rows(573).Resize(3).Insert()

# New row 573 - Crespo record / Primera
$ws.Cells.Item(573, 1).Value = 3
$ws.Cells.Item(573, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(573, 3).Value = "Coquimbo"
$ws.Cells.Item(573, 4).Value = 44706
$ws.Cells.Item(573, 5).Value = 5
$ws.Cells.Item(573, 6).Value = 100112006
$ws.Cells.Item(573, 7).Value = "Repollo"
$ws.Cells.Item(573, 8).Value = "Crespo record"
$ws.Cells.Item(573, 9).Value = "Primera"
$ws.Cells.Item(573, 10).Value = 2150
$ws.Cells.Item(573, 11).Value = 1000
$ws.Cells.Item(573, 12).Value = 1100
$ws.Cells.Item(573, 13).Value = 1056
$ws.Cells.Item(573, 14).Value = "$/unidad"
$ws.Cells.Item(573, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(573, 16).Value = 1056
$ws.Cells.Item(573, 17).Value = 1
$ws.Cells.Item(573, 18).Value = "Hortaliza"

# New row 574 - Crespo record / Segunda
$ws.Cells.Item(574, 1).Value = 3
$ws.Cells.Item(574, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(574, 3).Value = "Coquimbo"
$ws.Cells.Item(574, 4).Value = 44706
$ws.Cells.Item(574, 5).Value = 5
$ws.Cells.Item(574, 6).Value = 100112006
$ws.Cells.Item(574, 7).Value = "Repollo"
$ws.Cells.Item(574, 8).Value = "Crespo record"
$ws.Cells.Item(574, 9).Value = "Segunda"
$ws.Cells.Item(574, 10).Value = 1100
$ws.Cells.Item(574, 11).Value = 800
$ws.Cells.Item(574, 12).Value = 800
$ws.Cells.Item(574, 13).Value = 800
$ws.Cells.Item(574, 14).Value = "$/unidad"
$ws.Cells.Item(574, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(574, 16).Value = 800
$ws.Cells.Item(574, 17).Value = 1
$ws.Cells.Item(574, 18).Value = "Hortaliza"

# New row 575 - Morada(o) / Primera
$ws.Cells.Item(575, 1).Value = 3
$ws.Cells.Item(575, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(575, 3).Value = "Coquimbo"
$ws.Cells.Item(575, 4).Value = 44706
$ws.Cells.Item(575, 5).Value = 5
$ws.Cells.Item(575, 6).Value = 100112006
$ws.Cells.Item(575, 7).Value = "Repollo"
$ws.Cells.Item(575, 8).Value = "Morada(o)"
$ws.Cells.Item(575, 9).Value = "Primera"
$ws.Cells.Item(575, 10).Value = 1650
$ws.Cells.Item(575, 11).Value = 1500
$ws.Cells.Item(575, 12).Value = 1600
$ws.Cells.Item(575, 13).Value = 1548
$ws.Cells.Item(575, 14).Value = "$/unidad"
$ws.Cells.Item(575, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(575, 16).Value = 1548
$ws.Cells.Item(575, 17).Value = 1
$ws.Cells.Item(575, 18).Value = "Hortaliza"
